$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 3 (2019年) values per diff (rounding / value corrections) ---
$ws.Range("C3").Value = 9.699999999999999
$ws.Range("F3").Value = 3.4
$ws.Range("H3").Value = 50.5
$ws.Range("J3").Value = -1.2
$ws.Range("N3").Value = -15.6
$ws.Range("O3").Value = 0.3
$ws.Range("P3").Value = 2.9
$ws.Range("T3").Value = 1.4
$ws.Range("U3").Value = 0.7
$ws.Range("V3").Value = 1.6
$ws.Range("W3").Value = -8.699999999999999
$ws.Range("X3").Value = 3.1
$ws.Range("Y3").Value = 4.2
$ws.Range("Z3").Value = -14.1
$ws.Range("AA3").Value = 8.4
$ws.Range("AB3").Value = 5.3
$ws.Range("AC3").Value = 6.8
$ws.Range("AD3").Value = 4.6
$ws.Range("AM3").Value = -0.7
$ws.Range("AN3").Value = -9.1
$ws.Range("AQ3").Value = 21.9
$ws.Range("AR3").Value = -19.8
$ws.Range("AY3").Value = -15.9
$ws.Range("AZ3").Value = 17.7
$ws.Range("BA3").Value = 13.9
$ws.Range("BC3").Value = -2.4
$ws.Range("BE3").Value = 1.2
$ws.Range("BF3").Value = 6.8
$ws.Range("BI3").Value = -1.7
$ws.Range("BJ3").Value = 1
$ws.Range("BK3").Value = -22.5
$ws.Range("BL3").Value = 2.9
$ws.Range("BM3").Value = 1.4
$ws.Range("BN3").Value = 16.9
$ws.Range("BO3").Value = -1.5
$ws.Range("BP3").Value = 15.7
$ws.Range("BQ3").Value = -0.2
$ws.Range("BR3").Value = 29.6
$ws.Range("BS3").Value = 18.1
$ws.Range("BT3").Value = 37.2
$ws.Range("BU3").Value = 17.8
$ws.Range("BV3").Value = 4.5
$ws.Range("BW3").Value = -0.2
$ws.Range("BY3").Value = -3.6
$ws.Range("BZ3").Value = -2.6
$ws.Range("CA3").Value = 12.4
$ws.Range("CB3").Value = 25.7
$ws.Range("CF3").Value = 17.9
$ws.Range("CI3").Value = 15.8
$ws.Range("CJ3").Value = -3.4
$ws.Range("CK3").Value = -8.9
$ws.Range("CL3").Value = 1.8
$ws.Range("CN3").Value = -17.8
$ws.Range("CT3").Value = 2.2
$ws.Range("CU3").Value = -11.4
$ws.Range("CV3").Value = 9
$ws.Range("CX3").Value = 6.3
$ws.Range("CY3").Value = 24.1
$ws.Range("CZ3").Value = 38.1
$ws.Range("DA3").Value = -3.9
$ws.Range("DB3").Value = 10.4
$ws.Range("DC3").Value = -2.5
$ws.Range("DD3").Value = -0.1
$ws.Range("DF3").Value = 6.8
$ws.Range("DG3").Value = 30.9
$ws.Range("DH3").Value = -3.7
$ws.Range("DJ3").Value = 26
$ws.Range("DK3").Value = 2.5

# --- Add row 5 (2021年) ---
# Copy A4 (styled date cell) into A5 first to inherit its style (bold, centered, bordered),
# then overwrite with the correct label.
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("A5").Value = "2021年"
$ws.Range("B5").Value = 8
$ws.Range("C5").Value = 24.3
$ws.Range("D5").Value = -21.4
$ws.Range("E5").Value = 44.6
$ws.Range("F5").Value = 1.6
$ws.Range("G5").Value = 81.2
$ws.Range("H5").Value = 12
$ws.Range("I5").Value = 6.3
$ws.Range("J5").Value = 6.6
$ws.Range("K5").Value = -9.300000000000001
$ws.Range("L5").Value = -12.3
$ws.Range("M5").Value = -12.1
$ws.Range("N5").Value = -38.2
$ws.Range("O5").Value = -1.3
$ws.Range("P5").Value = 8.699999999999999
$ws.Range("Q5").Value = -12.3
$ws.Range("R5").Value = 9.199999999999999
$ws.Range("S5").Value = 18.1
$ws.Range("T5").Value = 10.5
$ws.Range("U5").Value = 9.300000000000001
$ws.Range("V5").Value = 3.4
$ws.Range("W5").Value = 18.8
$ws.Range("X5").Value = 13.5
$ws.Range("Y5").Value = 15.7
$ws.Range("Z5").Value = 31.8
$ws.Range("AA5").Value = 10.6
$ws.Range("AB5").Value = 19.5
$ws.Range("AC5").Value = 24.5
$ws.Range("AD5").Value = 6.7
$ws.Range("AE5").Value = 14
$ws.Range("AF5").Value = -40.8
$ws.Range("AH5").Value = -23.5
$ws.Range("AI5").Value = -7.7
$ws.Range("AJ5").Value = -16.9
$ws.Range("AK5").Value = 15.5
$ws.Range("AL5").Value = 4.7
$ws.Range("AM5").Value = 1.7
$ws.Range("AN5").Value = -10.3
$ws.Range("AO5").Value = -11.1
$ws.Range("AP5").Value = -27.1
$ws.Range("AQ5").Value = 6.7
$ws.Range("AR5").Value = 1.6
$ws.Range("AS5").Value = 43.7
$ws.Range("AT5").Value = 51.9
$ws.Range("AU5").Value = 45.2
$ws.Range("AV5").Value = 5
$ws.Range("AW5").Value = 10.8
$ws.Range("AX5").Value = -6.2
$ws.Range("AY5").Value = -5.9
$ws.Range("AZ5").Value = 11.7
$ws.Range("BA5").Value = 1.6
$ws.Range("BB5").Value = 6.4
$ws.Range("BC5").Value = 11.2
$ws.Range("BD5").Value = 58.5
$ws.Range("BE5").Value = 4.6
$ws.Range("BF5").Value = 1.9
$ws.Range("BG5").Value = 13.7
$ws.Range("BH5").Value = -3.3
$ws.Range("BI5").Value = -4
$ws.Range("BJ5").Value = 13.2
$ws.Range("BK5").Value = 17.9
$ws.Range("BL5").Value = -1.2
$ws.Range("BM5").Value = 1.3
$ws.Range("BN5").Value = -4.2
$ws.Range("BO5").Value = -3.7
$ws.Range("BP5").Value = 4.9
$ws.Range("BQ5").Value = 34.5
$ws.Range("BR5").Value = 11.1
$ws.Range("BS5").Value = -2.2
$ws.Range("BT5").Value = -2.6
$ws.Range("BU5").Value = -7
$ws.Range("BV5").Value = 1.1
$ws.Range("BW5").Value = 3.2
$ws.Range("BX5").Value = 23.3
$ws.Range("BY5").Value = 20.3
$ws.Range("BZ5").Value = 2
$ws.Range("CA5").Value = 8
$ws.Range("CB5").Value = 4.2
$ws.Range("CC5").Value = 17.3
$ws.Range("CD5").Value = -28.2
$ws.Range("CE5").Value = -3.1
$ws.Range("CF5").Value = 14.5
$ws.Range("CG5").Value = 15.8
$ws.Range("CH5").Value = 6.6
$ws.Range("CI5").Value = 13.6
$ws.Range("CJ5").Value = -1.4
$ws.Range("CK5").Value = 11.9
$ws.Range("CL5").Value = 4.1
$ws.Range("CM5").Value = -17.6
$ws.Range("CN5").Value = 18.8
$ws.Range("CO5").Value = 6.6
$ws.Range("CP5").Value = 22.3
$ws.Range("CQ5").Value = 10.6
$ws.Range("CR5").Value = -26.5
$ws.Range("CS5").Value = -15
$ws.Range("CT5").Value = 9.800000000000001
$ws.Range("CU5").Value = 13.3
$ws.Range("CV5").Value = -1.2
$ws.Range("CW5").Value = 25.2
$ws.Range("CX5").Value = 16.8
$ws.Range("CY5").Value = 10.9
$ws.Range("CZ5").Value = 43.6
$ws.Range("DA5").Value = 11.4
$ws.Range("DB5").Value = 1.9
$ws.Range("DC5").Value = 20.5
$ws.Range("DD5").Value = -1.8
$ws.Range("DE5").Value = -5.7
$ws.Range("DF5").Value = 14.1
$ws.Range("DG5").Value = 26.9
$ws.Range("DH5").Value = 10.4
$ws.Range("DI5").Value = 8.800000000000001
$ws.Range("DJ5").Value = 14.6
$ws.Range("DK5").Value = 26.9

# --- Add row 6 (2022年) ---
$ws.Range("A4").Copy($ws.Range("A6"))
$ws.Range("A6").Value = "2022年"
$ws.Range("C6").Value = 12.1
$ws.Range("F6").Value = 9.1
$ws.Range("H6").Value = 37.8
$ws.Range("J6").Value = 7.5
$ws.Range("M6").Value = 21.8
$ws.Range("N6").Value = 42.1
$ws.Range("O6").Value = 10.1
$ws.Range("P6").Value = -36.9
$ws.Range("T6").Value = 32.1
$ws.Range("U6").Value = 4.2
$ws.Range("V6").Value = 16.1
$ws.Range("W6").Value = 15.5
$ws.Range("X6").Value = 9.1
$ws.Range("Y6").Value = 18.8
$ws.Range("Z6").Value = 21.4
$ws.Range("AA6").Value = 5.9
$ws.Range("AB6").Value = 26.1
$ws.Range("AC6").Value = 27.3
$ws.Range("AD6").Value = 8
$ws.Range("AM6").Value = 13.2
$ws.Range("AN6").Value = 21.8
$ws.Range("AQ6").Value = 22.3
$ws.Range("AR6").Value = 2
$ws.Range("AV6").Value = -8.4
$ws.Range("AY6").Value = 5.3
$ws.Range("AZ6").Value = 5.4
$ws.Range("BA6").Value = 3.5
$ws.Range("BB6").Value = 9.4
$ws.Range("BC6").Value = 17.3
$ws.Range("BE6").Value = 15.7
$ws.Range("BF6").Value = 8.4
$ws.Range("BG6").Value = 19.6
$ws.Range("BI6").Value = 8.9
$ws.Range("BJ6").Value = 8.699999999999999
$ws.Range("BK6").Value = 16.5
$ws.Range("BL6").Value = 10.3
$ws.Range("BM6").Value = 13.6
$ws.Range("BN6").Value = 7.3
$ws.Range("BO6").Value = 12.6
$ws.Range("BP6").Value = 31.8
$ws.Range("BQ6").Value = -15
$ws.Range("BR6").Value = 24.4
$ws.Range("BS6").Value = 6.6
$ws.Range("BT6").Value = 6
$ws.Range("BU6").Value = 2.4
$ws.Range("BV6").Value = 19.3
$ws.Range("BW6").Value = 24.3
$ws.Range("BX6").Value = 42.6
$ws.Range("BY6").Value = -22.7
$ws.Range("BZ6").Value = 24.1
$ws.Range("CA6").Value = -10.7
$ws.Range("CB6").Value = 15.5
$ws.Range("CF6").Value = 21
$ws.Range("CI6").Value = 14.5
$ws.Range("CJ6").Value = 11.2
$ws.Range("CK6").Value = 4.7
$ws.Range("CL6").Value = 25.3
$ws.Range("CN6").Value = 4.8
$ws.Range("CP6").Value = 18.8
$ws.Range("CT6").Value = 14.8
$ws.Range("CU6").Value = 8.300000000000001
$ws.Range("CV6").Value = 3.7
$ws.Range("CX6").Value = 27.2
$ws.Range("CY6").Value = 4.5
$ws.Range("CZ6").Value = -22.4
$ws.Range("DA6").Value = 11.8
$ws.Range("DB6").Value = 10.5
$ws.Range("DC6").Value = 1.7
$ws.Range("DD6").Value = 1.8
$ws.Range("DF6").Value = 6.7
$ws.Range("DG6").Value = 17.3
$ws.Range("DH6").Value = 13.7
$ws.Range("DJ6").Value = -0.1
$ws.Range("DK6").Value = 33.3
